# This workbook holds weekly price records (one row per observation) for
# "Papa" (potato) at "Terminal Hortofrutícola Agro Chillán". A new weekly
# observation was inserted as row 151, pushing all subsequent rows down by
# one (old row 151 -> new row 152, ..., old row 195 -> new row 196).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 151; this shifts rows 151:195
# down to 152:196 and extends the used range to A1:R196.
$ws.Rows(151).Insert()

# Populate the newly inserted row 151 with the new weekly record.
$ws.Range("A151").Value = 7
$ws.Range("B151").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C151").Value = "Ñuble"
$ws.Range("D151").Value = 44463
$ws.Range("E151").Value = 16
$ws.Range("F151").Value = 100114001
$ws.Range("G151").Value = "Papa"
$ws.Range("H151").Value = "Patagonia"
$ws.Range("I151").Value = "1a (guarda)"
$ws.Range("J151").Value = 300
$ws.Range("K151").Value = 6500
$ws.Range("L151").Value = 7000
$ws.Range("M151").Value = 6750
$ws.Range("N151").Value = '$/saco 25 kilos'
$ws.Range("O151").Value = "Provincia de Diguillín"
$ws.Range("P151").Value = 270
$ws.Range("Q151").Value = 25
$ws.Range("R151").Value = "Hortaliza"
